$d = $word.ActiveDocument

$d.Content.Find.Execute("20÷7=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=7, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷3=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=11, 1", 2) | Out-Null
$d.Content.Find.Execute("33÷6=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2) | Out-Null
$d.Content.Find.Execute("29÷4=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "43÷4=10, 3", 2) | Out-Null
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=7, 5", 2) | Out-Null
$d.Content.Find.Execute("69÷8=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2) | Out-Null
$d.Content.Find.Execute("11÷6=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=2, 0", 2) | Out-Null
$d.Content.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷8=5, 6", 2) | Out-Null
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "28÷3=9, 1", 2) | Out-Null
$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "80÷7=11, 3", 2) | Out-Null
$d.Content.Find.Execute("34÷9=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "67÷4=16, 3", 2) | Out-Null
$d.Content.Find.Execute("66÷3=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "37÷2=18, 1", 2) | Out-Null
$d.Content.Find.Execute("32÷7=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "82÷9=9, 1", 2) | Out-Null
$d.Content.Find.Execute("94÷7=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2) | Out-Null
$d.Content.Find.Execute("98÷9=10, 8", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=5, 0", 2) | Out-Null
$d.Content.Find.Execute("83÷2=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "27÷4=6, 3", 2) | Out-Null
$d.Content.Find.Execute("40÷9=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "83÷4=20, 3", 2) | Out-Null
$d.Content.Find.Execute("38÷2=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "38÷4=9, 2", 2) | Out-Null
$d.Content.Find.Execute("48÷7=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "64÷5=12, 4", 2) | Out-Null
$d.Content.Find.Execute("54÷3=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=8, 4", 2) | Out-Null
$d.Content.Find.Execute("27÷5=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=7, 5", 2) | Out-Null
$d.Content.Find.Execute("93÷6=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "37÷6=6, 1", 2) | Out-Null
$d.Content.Find.Execute("58÷5=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=6, 5", 2) | Out-Null
$d.Content.Find.Execute("70÷2=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷4=24, 1", 2) | Out-Null
